$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For Price cells whose new value is a plain decimal number (e.g. "608.35"),
# force the cell to keep Text formatting first so Excel does not silently
# reinterpret the assigned string as a floating point number (these columns
# store formatted price strings, not numeric values).

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.980.89'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '3.539.31'
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '608.35'
$ws.Range("E5").Value = '  +3.03%  '
$ws.Range("D6").Value = '185.54'
$ws.Range("E6").Value = '  -1.17%  '
$ws.Range("D7").Value = '3.535.64'
$ws.Range("E7").Value = '  -1.02%  '
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").Value = '0.214'
$ws.Range("E10").Value = '  +5.56%  '
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("D12").Value = '53.69'
$ws.Range("E12").Value = '  -2.33%  '
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("D14").Value = '9.46'
$ws.Range("E14").Value = '  -1.85%  '
$ws.Range("D15").Value = '4.098.03'
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D16").Value = '70.056.58'
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("D17").Value = '12.65'
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.548.43'
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '18.89'
$ws.Range("E19").Value = '  -3.15%  '
$ws.Range("D20").Value = '577.22'
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").Value = '0.991'
$ws.Range("E22").Value = '  -3.26%  '
$ws.Range("D23").Value = '17.44'
$ws.Range("E23").Value = '  -2.42%  '
$ws.Range("D24").Value = '4.67'
$ws.Range("E24").Value = '  -0.84%  '
$ws.Range("D25").Value = '4.88'
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("D26").Value = '95.26'
$ws.Range("E26").Value = '  -1.06%  '
$ws.Range("D27").Value = '2.98'
$ws.Range("E27").Value = '  -1.06%  '
$ws.Range("D28").Value = '10.96'
$ws.Range("E28").Value = '  -4.94%  '
$ws.Range("D29").Value = '9.41'
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("E30").Value = '  -1.25%  '
$ws.Range("E31").Value = '  -4.85%  '
$ws.Range("D32").Value = '12.17'
$ws.Range("E32").Value = '  -3.40%  '
$ws.Range("E33").Value = '  -1.91%  '
$ws.Range("D34").Value = '63.09'
$ws.Range("E34").Value = '  -3.23%  '
$ws.Range("B35").Value = 'dogwifhat'
$ws.Range("C35").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D35").Value = '3.62'
$ws.Range("E35").Value = '  +16.43%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").Value = '3.26'
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("D37").Value = '536.79'
$ws.Range("E37").Value = '  -5.32%  '
$ws.Range("E38").Value = '  -3.89%  '
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").Value = '37.13'
$ws.Range("E40").Value = '  -3.19%  '
$ws.Range("D41").Value = '0.0₃0778'
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("D42").Value = '3.534.28'
$ws.Range("D43").Value = '3.53'
$ws.Range("E43").Value = '  +3.93%  '
$ws.Range("E44").Value = '  +0.38%  '
$ws.Range("E45").Value = '  +1.31%  '
$ws.Range("E46").Value = '  -2.36%  '
$ws.Range("D47").Value = '3.37'
$ws.Range("E47").Value = '  -5.26%  '
$ws.Range("E48").Value = '  +2.41%  '
$ws.Range("D49").Value = '9.12'
$ws.Range("E49").Value = '  -3.99%  '
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("D51").Value = '1.41'
$ws.Range("E51").Value = '  -4.37%  '
